$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 398.8
$ws.Range("I70").Value = 399.79797
$ws.Range("J70").Value = 300
$ws.Range("K70").Value = 1199.39391
$ws.Range("L70").Value = 900
$ws.Range("M70").Value = -929.39391
$ws.Range("N70").Value = -1440

$ws.Range("H73").Value = 398.8
$ws.Range("I73").Value = 399.79797
$ws.Range("J73").Value = 300
$ws.Range("K73").Value = 1199.39391
$ws.Range("L73").Value = 900
$ws.Range("M73").Value = -263.39391
$ws.Range("N73").Value = -2772

$ws.Range("H100").Value = 6187.5
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 8100
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 8100
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -9182

$ws.Range("H133").Value = 48799.668
$ws.Range("J133").Value = 48799.668
$ws.Range("L133").Value = 48799.668
$ws.Range("N133").Value = -58919.668

$ws.Range("H134").Value = 49099.668
$ws.Range("J134").Value = 49099.668
$ws.Range("L134").Value = 49099.668
$ws.Range("N134").Value = -59239.668

$ws.Range("H136").Value = 48899.8
$ws.Range("J136").Value = 48899.8
$ws.Range("L136").Value = 48899.8
$ws.Range("N136").Value = -59099.8

$ws.Range("H139").Value = 45664.75
$ws.Range("J139").Value = 45664.75
$ws.Range("L139").Value = 45664.75
$ws.Range("N139").Value = -55944.75

$ws.Range("H140").Value = 37294.75
$ws.Range("J140").Value = 37294.75
$ws.Range("L140").Value = 37294.75
$ws.Range("N140").Value = -47654.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1998.6471
$ws.Range("I45").Value = 1997.0769
$ws.Range("K45").Value = 1997.0769
$ws.Range("M45").Value = -1620.0769

$ws.Range("H61").Value = 1916.0555
$ws.Range("I61").Value = 1414.0952
$ws.Range("J61").Value = 2618.8
$ws.Range("K61").Value = 1414.0952
$ws.Range("L61").Value = 2618.8
$ws.Range("M61").Value = -1202.0952
$ws.Range("N61").Value = -3042.8

$ws.Range("H63").Value = 3398.889
$ws.Range("I63").Value = 2427.1428
$ws.Range("K63").Value = 2427.1428
$ws.Range("M63").Value = -1741.1428

$ws.Range("H66").Value = 3398.889
$ws.Range("I66").Value = 2427.1428
$ws.Range("K66").Value = 12135.714
$ws.Range("M66").Value = -8703.714

$ws.Range("H74").Value = 9204727
$ws.Range("I74").Value = 7258511
$ws.Range("J74").Value = 33337800
$ws.Range("K74").Value = 7258511
$ws.Range("L74").Value = 33337800
$ws.Range("M74").Value = -7257637
$ws.Range("N74").Value = -33339548

$ws.Range("H77").Value = 9204727
$ws.Range("I77").Value = 7258511
$ws.Range("J77").Value = 33337800
$ws.Range("K77").Value = 36292555
$ws.Range("L77").Value = 166689000
$ws.Range("M77").Value = -36288187
$ws.Range("N77").Value = -166697736

$ws.Range("H136").Value = 1916.0555
$ws.Range("I136").Value = 1414.0952
$ws.Range("J136").Value = 2618.8
$ws.Range("K136").Value = 4242.2856
$ws.Range("L136").Value = 7856.400000000001
$ws.Range("M136").Value = -1692.2856
$ws.Range("N136").Value = -12956.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 105000
$ws.Range("J42").Value = 105000
$ws.Range("L42").Value = 105000
$ws.Range("N42").Value = -105656

$ws.Range("H86").Value = 22226222
$ws.Range("J86").Value = 9000
$ws.Range("L86").Value = 9000
$ws.Range("N86").Value = -11246

$ws.Range("H89").Value = 22226222
$ws.Range("J89").Value = 9000
$ws.Range("L89").Value = 45000
$ws.Range("N89").Value = -56232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 25000
$ws.Range("J32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("N32").Value = -25632

$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16498

$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -52488

$ws.Range("H86").Value = 251616.84
$ws.Range("J86").Value = 1836.6666
$ws.Range("L86").Value = 1836.6666
$ws.Range("N86").Value = -4082.6666

$ws.Range("H89").Value = 251616.84
$ws.Range("J89").Value = 1836.6666
$ws.Range("L89").Value = 9183.333000000001
$ws.Range("N89").Value = -20415.333

$ws.Range("H99").Value = 8940149
$ws.Range("I99").Value = 11915865
$ws.Range("K99").Value = 11915865
$ws.Range("M99").Value = -11914367

$ws.Range("H126").Value = 8940149
$ws.Range("I126").Value = 11915865
$ws.Range("K126").Value = 35747595
$ws.Range("M126").Value = -35745125

$ws.Range("H134").Value = 2968.484
$ws.Range("I134").Value = 3182.682
$ws.Range("J134").Value = 2444.889
$ws.Range("K134").Value = 9548.045999999998
$ws.Range("L134").Value = 7334.667
$ws.Range("M134").Value = -7013.045999999998
$ws.Range("N134").Value = -12404.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 9629.625
$ws.Range("I63").Value = 2184
$ws.Range("J63").Value = 17075.25
$ws.Range("K63").Value = 6552
$ws.Range("L63").Value = 51225.75
$ws.Range("M63").Value = -5803
$ws.Range("N63").Value = -52723.75

$ws.Range("H66").Value = 9629.625
$ws.Range("I66").Value = 2184
$ws.Range("J66").Value = 17075.25
$ws.Range("K66").Value = 19656
$ws.Range("L66").Value = 153677.25
$ws.Range("M66").Value = -15912
$ws.Range("N66").Value = -161165.25

$ws.Range("H107").Value = 1151.8422
$ws.Range("I107").Value = 317.33334
$ws.Range("J107").Value = 1410.8276
$ws.Range("K107").Value = 952.0000200000001
$ws.Range("L107").Value = 4232.4828
$ws.Range("M107").Value = 967.9999799999999
$ws.Range("N107").Value = -8072.4828

$ws.Range("H112").Value = 7771.3076
$ws.Range("I112").Value = 2750.8
$ws.Range("J112").Value = 8966.666999999999
$ws.Range("K112").Value = 8252.400000000001
$ws.Range("L112").Value = 26900.001
$ws.Range("M112").Value = -7144.400000000001
$ws.Range("N112").Value = -29116.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3803.35
$ws.Range("J80").Value = 2934.3333
$ws.Range("L80").Value = 2934.3333
$ws.Range("N80").Value = -4930.3333

$ws.Range("H83").Value = 3803.35
$ws.Range("J83").Value = 2934.3333
$ws.Range("L83").Value = 14671.6665
$ws.Range("N83").Value = -24655.6665

$ws.Range("H132").Value = 2713.1304
$ws.Range("I132").Value = 1938.7693
$ws.Range("J132").Value = 3719.8
$ws.Range("K132").Value = 5816.3079
$ws.Range("L132").Value = 11159.4
$ws.Range("M132").Value = -3286.3079
$ws.Range("N132").Value = -16219.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4154.6665
$ws.Range("I7").Value = 4004
$ws.Range("J7").Value = 4230
$ws.Range("K7").Value = 4004
$ws.Range("L7").Value = 4230
$ws.Range("M7").Value = -3892
$ws.Range("N7").Value = -4454

$ws.Range("H16").Value = 1473
$ws.Range("I16").Value = 1559.091
$ws.Range("K16").Value = 1559.091
$ws.Range("M16").Value = -1389.091

$ws.Range("H40").Value = 4652.8696
$ws.Range("I40").Value = 4579.8423
$ws.Range("J40").Value = 4999.75
$ws.Range("K40").Value = 4579.8423
$ws.Range("L40").Value = 4999.75
$ws.Range("M40").Value = -4443.8423
$ws.Range("N40").Value = -5271.75

$ws.Range("H55").Value = 239.64706
$ws.Range("I55").Value = 187
$ws.Range("J55").Value = 314.85715
$ws.Range("K55").Value = 187
$ws.Range("L55").Value = 314.85715
$ws.Range("M55").Value = -14
$ws.Range("N55").Value = -660.85715

$ws.Range("H100").Value = 1271
$ws.Range("I100").Value = 1052.3
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1052.3
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -511.3
$ws.Range("N100").Value = -3082

$ws.Range("H126").Value = 4154.6665
$ws.Range("I126").Value = 4004
$ws.Range("J126").Value = 4230
$ws.Range("K126").Value = 12012
$ws.Range("L126").Value = 12690
$ws.Range("M126").Value = -9542
$ws.Range("N126").Value = -17630

$ws.Range("H136").Value = 2843871.2
$ws.Range("I136").Value = 4170465.5
$ws.Range("J136").Value = 1169.2142
$ws.Range("K136").Value = 12511396.5
$ws.Range("L136").Value = 3507.6426
$ws.Range("M136").Value = -12508846.5
$ws.Range("N136").Value = -8607.642599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 605520
$ws.Range("J62").Value = 1005200
$ws.Range("L62").Value = 1005200
$ws.Range("N62").Value = -1006448

$ws.Range("H65").Value = 605520
$ws.Range("J65").Value = 1005200
$ws.Range("L65").Value = 5026000
$ws.Range("N65").Value = -5032240

$ws.Range("H113").Value = 100001760
$ws.Range("I113").Value = 71430800
$ws.Range("J113").Value = 166667330
$ws.Range("K113").Value = 214292400
$ws.Range("L113").Value = 500001990
$ws.Range("M113").Value = -214290230
$ws.Range("N113").Value = -500006330

$ws.Range("H122").Value = 142860240
$ws.Range("I122").Value = 200003140
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 600009420
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -600006970
$ws.Range("N122").Value = -13900

$ws.Range("H133").Value = 43500
$ws.Range("J133").Value = 43500
$ws.Range("L133").Value = 43500
$ws.Range("N133").Value = -53620

$ws.Range("H136").Value = 4449.9355
$ws.Range("I136").Value = 623
$ws.Range("J136").Value = 24350
$ws.Range("K136").Value = 1869
$ws.Range("L136").Value = 73050
$ws.Range("M136").Value = 681
$ws.Range("N136").Value = -78150
